$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the second data row (old row 3), leaving only header + 1 data row
$ws.Rows.Item(3).Delete()

# Add new header columns L and M, matching the bold header style used in A1:K1
$ws.Range("L1").Value = "Criador do Registro"
$ws.Range("M1").Value = "Validador do Registro"
$ws.Range("A1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)

# Force the data row to be stored as text (avoids numeric/date auto-conversion of
# numeric-looking or date-looking strings), then strip the resulting formatting
# so cells keep the workbook's default (unstyled) look.
$dataRow = $ws.Range("A2:M2")
$dataRow.NumberFormat = "@"

$ws.Range("A2").Value = "CC_1609929113625"
$ws.Range("B2").Value = "Obra_001"
$ws.Range("C2").Value = "OAK-4403"
$ws.Range("D2").Value = "BRITA 1"
$ws.Range("E2").Value = "341"
$ws.Range("F2").Value = "841"
$ws.Range("G2").Value = "-3.07022208392808"
$ws.Range("H2").Value = "-60.0082966808251"
$ws.Range("I2").Value = "6-1-2021"
$ws.Range("J2").Value = "6:31"
$ws.Range("K2").Value = "6:32"
$ws.Range("L2").Value = "Admin"
$ws.Range("M2").Value = "Admin"

$dataRow.ClearFormats()

# Set column widths for the newly added columns (existing cols 1-11 already at width 20):
# col 12 (L) = 20, col 13 (M) = 21. COM ColumnWidth is in "characters" and gets
# padding added internally, so nudge the input slightly above the integer boundary
# to land exactly on the target stored width.
$ws.Columns.Item(12).ColumnWidth = 19.1
$ws.Columns.Item(13).ColumnWidth = 20.1
